# Update countries & provincias Spain
#
# Refreshes the "Pais" worksheet (daily COVID-19 stats) with the latest
# figures. Several countries changed rank when the source table was
# re-sorted by "Casos totales", so a handful of rows need both their
# country label (column A) and their figures (columns B:H) replaced;
# the rest of the changed rows only need updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 18:29"

# --- Estados Unidos ---
$ws.Range("B4").Value = 4781939
$ws.Range("C4").Value = 17621
$ws.Range("D4").Value = 2365246
$ws.Range("E4").Value = 2258612
$ws.Range("G4").Value = 183
$ws.Range("H4").Value = 158081

# --- India ---
$ws.Range("B6").Value = 1796486
$ws.Range("C6").Value = 44567
$ws.Range("D6").Value = 1181134
$ws.Range("E6").Value = 577289
$ws.Range("G6").Value = 660
$ws.Range("H6").Value = 38063

# --- Reino Unido ---
$ws.Range("B15").Value = 304695
$ws.Range("C15").Value = 743
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 46201

# --- Italia ---
$ws.Range("B18").Value = 248070
$ws.Range("C18").Value = 238
$ws.Range("D18").Value = 200460
$ws.Range("E18").Value = 12456
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 35154

# --- Alemania ---
$ws.Range("B21").Value = 211257
$ws.Range("C21").Value = 180
$ws.Range("E21").Value = 8431

# --- Canada ---
$ws.Range("B25").Value = 116858
$ws.Range("C25").Value = 259
$ws.Range("D25").Value = 101558
$ws.Range("E25").Value = 6355
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 8945

# --- Row 37 now ranks "Republica Dominicana" (was "Ucrania") ---
$ws.Range("A37").Value = 'Republica Dominicana'
$ws.Range("B37").Value = 72243
$ws.Range("C37").Value = 828
$ws.Range("D37").Value = 38244
$ws.Range("E37").Value = 32821
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 1178

# --- Row 38 now ranks "Ucrania" (was "Republica Dominicana") ---
$ws.Range("A38").Value = 'Ucrania'
$ws.Range("B38").Value = 72168
$ws.Range("C38").Value = 1112
$ws.Range("D38").Value = 39543
$ws.Range("E38").Value = 30900
$ws.Range("G38").Value = 16
$ws.Range("H38").Value = 1725

# --- Guatemala ---
$ws.Range("B48").Value = 51306
$ws.Range("C48").Value = 327
$ws.Range("D48").Value = 38416
$ws.Range("E48").Value = 10895
$ws.Range("G48").Value = 36
$ws.Range("H48").Value = 1995

# --- Azerbaiyan ---
$ws.Range("B59").Value = 32443
$ws.Range("C59").Value = 286
$ws.Range("D59").Value = 27113
$ws.Range("E59").Value = 4868
$ws.Range("G59").Value = 8
$ws.Range("H59").Value = 462

# --- Argelia ---
$ws.Range("B60").Value = 31465
$ws.Range("C60").Value = 515
$ws.Range("D60").Value = 21419
$ws.Range("E60").Value = 8815
$ws.Range("G60").Value = 8
$ws.Range("H60").Value = 1231

# --- Kenia ---
$ws.Range("D66").Value = 8477
$ws.Range("E66").Value = 13207

# --- Chequia ---
$ws.Range("B75").Value = 16729
$ws.Range("C75").Value = 30
$ws.Range("D75").Value = 11596
$ws.Range("E75").Value = 4749
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 384

# --- Row 104 now ranks "Grecia" (was "Republica de Africa Central") ---
$ws.Range("A104").Value = 'Grecia'
$ws.Range("B104").Value = 4662
$ws.Range("C104").Value = 75
$ws.Range("D104").Value = 1374
$ws.Range("E104").Value = 3080
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 208

# --- Row 105 now ranks "Republica de Africa Central" (was "Grecia") ---
$ws.Range("A105").Value = 'Republica de Africa Central'
$ws.Range("B105").Value = 4614
$ws.Range("D105").Value = 1635
$ws.Range("E105").Value = 2920
$ws.Range("H105").Value = 59

# --- Sri Lanka ---
$ws.Range("B118").Value = 2822
$ws.Range("C118").Value = 7
$ws.Range("E118").Value = 297

# --- Estonia ---
$ws.Range("E128").Value = 82
$ws.Range("H128").Value = 63

# --- Liberia ---
$ws.Range("B142").Value = 1207
$ws.Range("C142").Value = 18
$ws.Range("E142").Value = 457
$ws.Range("G142").Value = 2
$ws.Range("H142").Value = 77

# --- Row 146 now ranks "Republica de Chipre" (was "Burkina Faso") ---
$ws.Range("A146").Value = 'Republica de Chipre'
$ws.Range("B146").Value = 1150
$ws.Range("C146").Value = 26
$ws.Range("D146").Value = 856
$ws.Range("E146").Value = 275
$ws.Range("H146").Value = 19

# --- Row 147 now ranks "Burkina Faso" (was "Republica de Chipre") ---
$ws.Range("A147").Value = 'Burkina Faso'
$ws.Range("B147").Value = 1143
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 945
$ws.Range("E147").Value = 145
$ws.Range("H147").Value = 53
